# The sheet contains a weekly price series for "Perejil" (parsley) at the
# Vega Central Mapocho de Santiago market. A new weekly record needs to be
# inserted as row 192 (pushing every following record down by one row), so
# the table grows from 240 rows (A1:R240) to 241 rows (A1:R241).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new blank row at position 192; everything that was in rows
# 192-240 shifts down to 193-241 (Excel keeps the existing data/format intact).
$ws.Rows("192:192").Insert()

# Populate the newly inserted row 192 with the new weekly record.
$ws.Cells.Item(192, 1).Value  = 9
$ws.Cells.Item(192, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(192, 3).Value  = "Metropolitana"
$ws.Cells.Item(192, 4).Value  = 44511
$ws.Cells.Item(192, 5).Value  = 13
$ws.Cells.Item(192, 6).Value  = 100112044
$ws.Cells.Item(192, 7).Value  = "Perejil"
$ws.Cells.Item(192, 8).Value  = "Sin especificar"
$ws.Cells.Item(192, 9).Value  = "Primera"
$ws.Cells.Item(192, 10).Value = 79
$ws.Cells.Item(192, 11).Value = 14000
$ws.Cells.Item(192, 12).Value = 15000
$ws.Cells.Item(192, 13).Value = 14506
$ws.Cells.Item(192, 14).Value = "$/docena de atados"
$ws.Cells.Item(192, 15).Value = "Región Metropolitana"
$ws.Cells.Item(192, 16).Value = 4835
$ws.Cells.Item(192, 17).Value = 3
$ws.Cells.Item(192, 18).Value = "Hortaliza"
